$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume(1h) (E) columns keep their values as plain text,
# matching how the data was originally stored (inline/shared strings), so values
# like "0.110" or "2.70" are not silently reinterpreted as numbers and lose formatting.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '51.879.71'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '2.992.23'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '354.68'
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").Value = '106.49'
$ws.Range("E6").Value = '  -4.07%  '
$ws.Range("D7").Value = '0.554'
$ws.Range("E7").Value = '  -2.39%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").Value = '0.605'
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").Value = '37.84'
$ws.Range("E10").Value = '  -4.20%  '
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").Value = '0.0853'
$ws.Range("E12").Value = '  -3.24%  '
$ws.Range("E13").Value = '  -3.61%  '
$ws.Range("D14").Value = '3.472.39'
$ws.Range("E14").Value = '  +2.37%  '
$ws.Range("D15").Value = '7.55'
$ws.Range("E15").Value = '  -4.51%  '
$ws.Range("D16").Value = '2.983.98'
$ws.Range("D17").Value = '0.998'
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D18").Value = '51.907.53'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = '3.36'
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").Value = '7.40'
$ws.Range("E20").Value = '  -2.34%  '
$ws.Range("D21").Value = '13.42'
$ws.Range("E21").Value = '  -4.61%  '
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("D23").Value = '68.90'
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("D24").Value = '262.93'
$ws.Range("E24").Value = '  -2.90%  '
$ws.Range("D25").Value = '2.70'
$ws.Range("E25").Value = '  -4.25%  '
$ws.Range("D26").Value = '0.177'
$ws.Range("E26").Value = '  -4.23%  '
$ws.Range("D27").Value = '26.78'
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").Value = '0.110'
$ws.Range("E30").Value = '  +2.86%  '
$ws.Range("D31").Value = '6.29'
$ws.Range("E31").Value = '  +4.01%  '
$ws.Range("D32").Value = '10.11'
$ws.Range("E32").Value = '  -4.68%  '
$ws.Range("D33").Value = '35.82'
$ws.Range("E33").Value = '  -7.62%  '
$ws.Range("D34").Value = '2.16'
$ws.Range("E34").Value = '  +12.48%  '
$ws.Range("D35").Value = '50.98'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").Value = '0.0429'
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").Value = '2.81'
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("D40").Value = '1.93'
$ws.Range("E40").Value = '  -3.90%  '
$ws.Range("D41").Value = '17.42'
$ws.Range("E41").Value = '  -6.14%  '
$ws.Range("E42").Value = '  -3.42%  '
$ws.Range("D43").Value = '23.13'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '124.37'
$ws.Range("E44").Value = '  +4.20%  '
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").Value = '2.115.67'
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("D47").Value = '3.29'
$ws.Range("E47").Value = '  -5.14%  '
$ws.Range("D49").Value = '3.298.72'
$ws.Range("E49").Value = '  +2.46%  '
$ws.Range("D50").Value = '0.241'
$ws.Range("E50").Value = '  -3.17%  '
$ws.Range("D51").Value = '0.0330'
$ws.Range("E51").Value = '  -0.93%  '
